# Update CA Excel file with latest revisions:
#  - Fill in the "Written Assignment (15)" (column D) scores that were
#    previously recorded as 0.
#  - Remove the derived "Total (60)" column (I) - header, formulas and
#    the now-unused shared string - since totals are no longer tracked
#    in this sheet.
#  - Leave the UI selection state close to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column D (Written Assignment) values, rows 2-27.
$dValues = @{
    2  = 10.8
    3  = 10.3
    4  = 10.5
    5  = 11.3
    6  = 11.5
    7  = 11.5
    8  = 10.5
    9  = 10.5
    10 = 10.5
    11 = 10.5
    12 = 10.3
    13 = 11.5
    14 = 10.3
    15 = 10.3
    16 = 10.5
    17 = 11.3
    18 = 10.8
    19 = 10.8
    20 = 11.5
    21 = 10.8
    22 = 10.8
    23 = 10.3
    24 = 11.3
    25 = 10.3
    26 = 11.5
    27 = 10.5
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

# 2. Mark up the selection the way the workbook was left (best-effort -
#    the second, active, member of the multi-area selection).
$ws.Range("D4").Select() | Out-Null
$ws.Range("I16").Activate() | Out-Null

# 3. Remove the whole "Total (60)" column (I) - header text, the
#    shared-formula totals and the trailing shared string all go away
#    with it, and the used range shrinks back to A1:H27.
$ws.Columns.Item(9).Delete() | Out-Null
